$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - column headers
$ws.Range("A1").Value = "Comarca nombre"
$ws.Range("B1").Value = "Número hogares"
$ws.Range("C1").Value = "Comarca código"
$ws.Range("D1").Value = "Provincia código"
$ws.Range("E1").Value = "Aragón"
$ws.Range("F1").Value = "Municipio código"
$ws.Range("G1").Value = "Provincia nombre"
$ws.Range("H1").Value = "Año de construcción"
$ws.Range("I1").Value = "Municipio nombre"

# Row 2 - sdmx / iaest dimension identifiers
$ws.Range("A2").Value = "sdmx-dimension:refArea"
$ws.Range("B2").Value = "iaest-measure:numero-hogares"
$ws.Range("D2").Value = "null"
$ws.Range("F2").Value = "null"
$ws.Range("G2").Value = "sdmx-dimension:refArea"
$ws.Range("H2").Value = "iaest-dimension:ano-de-construccion"
$ws.Range("I2").Value = "sdmx-dimension:refArea"

# Row 3 - dim/medida markers
$ws.Range("A3").Value = "dim"
$ws.Range("B3").Value = "medida"
$ws.Range("D3").Value = "null"
$ws.Range("F3").Value = "null"
$ws.Range("G3").Value = "dim"
$ws.Range("I3").Value = "dim"

# Row 4 - URI / type values
$ws.Range("A4").Value = "URI-comarca"
$ws.Range("B4").Value = "xsd:int"
$ws.Range("D4").Value = "null"
$ws.Range("E4").Value = "URI-Comunidad"
$ws.Range("F4").Value = "null"
$ws.Range("G4").Value = "URI-Provincia"
$ws.Range("H4").Value = "skos:Concept"
$ws.Range("I4").Value = "URI-Municipio"

# Row 5 - mapping file moved from B5 to H5
$ws.Range("H5").Value = "mapping-ano-de-construccion.xlsx"
$ws.Range("B4").Copy()
$ws.Range("H5").PasteSpecial(-4122)
$ws.Range("B5").Clear()
